$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315593004226685
$ws.Range("B1").Value = 1.452528476715088
$ws.Range("C1").Value = 4.493249416351318
$ws.Range("D1").Value = 5.08116340637207
$ws.Range("E1").Value = 1.517598748207092
